$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.817.41"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.884.25"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'0.7451"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.61%  "

# Row 6
$ws.Range("D6").Value = "'241.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.3110"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.89%  "

# Row 9
$ws.Range("D9").Value = "'25.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.55%  "

# Row 10
$ws.Range("D10").Value = "'0.07076"
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.08488"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.12%  "

# Row 12
$ws.Range("D12").Value = "'0.7563"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.04%  "

# Row 13
$ws.Range("D13").Value = "'1.892.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "

# Row 14
$ws.Range("D14").Value = "'5.349"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.89%  "

# Row 15
$ws.Range("D15").Value = "'92.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.04%  "

# Row 16
$ws.Range("D16").Value = "'6.121"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.66%  "

# Row 17
$ws.Range("D17").Value = "'29.816.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.44%  "

# Row 18
$ws.Range("D18").Value = "'13.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.58%  "

# Row 19
$ws.Range("D19").Value = "'242.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.95%  "

# Row 20
$ws.Range("D20").Value = "'0.000007800"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.71%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.142.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "

# Row 23
$ws.Range("D23").Value = "'7.965"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.36%  "

# Row 24
$ws.Range("D24").Value = "'1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").Value = "'0.1575"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "

# Row 26
$ws.Range("D26").Value = "'9.316"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.38%  "

# Row 27
$ws.Range("D27").Value = "'162.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "

# Row 28
$ws.Range("D28").Value = "'18.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.68%  "

# Row 29
$ws.Range("D29").Value = "'2.021"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.43%  "

# Row 30
$ws.Range("D30").Value = "'1.474"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.51%  "

# Row 31
$ws.Range("D31").Value = "'1.530"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

# Row 32
$ws.Range("D32").Value = "'4.480"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.11%  "

# Row 33
$ws.Range("D33").Value = "'4.151"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.05%  "

# Row 34
$ws.Range("D34").Value = "'0.05386"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.53%  "

# Row 35
$ws.Range("D35").Value = "'1.234"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.42%  "

# Row 36
$ws.Range("D36").Value = "'0.7502"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "

# Row 37
$ws.Range("D37").Value = "'1.007"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "

# Row 38
$ws.Range("D38").Value = "'2.703"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "

# Row 39
$ws.Range("D39").Value = "'0.01935"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "

# Row 40
$ws.Range("D40").Value = "'2.769"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.80%  "

# Row 41
$ws.Range("D41").Value = "'0.4449"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.088"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.24%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'1.095.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.67%  "

# Row 44
$ws.Range("D44").Value = "'72.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.08%  "

# Row 45
$ws.Range("D45").Value = "'0.8660"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.75%  "

# Row 46
$ws.Range("D46").Value = "'1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'102.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.12%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.690"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.93%  "

# Row 49
$ws.Range("D49").Value = "'1.847"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.33%  "

# Row 50
$ws.Range("D50").Value = "'3.017"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.66%  "

# Row 51
$ws.Range("D51").Value = "'2.034.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.63%  "
